$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a brand-new task paragraph right after "플레이어 움직임 구현"
#    (and therefore right before "히트 체크 구현"), while neither of
#    its neighbours carries strikethrough yet so nothing bleeds into
#    the new paragraph's inherited paragraph-mark formatting:
#       "       6-2. 플레이어 crouch(앉기) 구현"
#    with a bookmark wrapped around "crouch".
# ------------------------------------------------------------------
$tail = $d.Paragraphs.Item(8).Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$newPar = $d.Paragraphs.Item(9)
$newPar.Range.ListFormat.RemoveNumbers()
$newPar.Style = "List Paragraph"
$newPar.Range.ParagraphFormat.LeftIndent = 40

$pos = $newPar.Range.Start

function Add-Piece($text, $color) {
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($text)
    $r2 = $d.Range($pos, $pos + $text.Length)
    $r2.Font.Color = $color
    $pos = $pos + $text.Length
    return $r2
}

$null       = Add-Piece "       " 255
$null       = Add-Piece "6-2. " 0
$null       = Add-Piece "플레이어 " 0
$crouchRng  = Add-Piece "crouch" 0
$null       = Add-Piece "(" 0
$null       = Add-Piece "앉기" 0
$null       = Add-Piece ") " 0
$null       = Add-Piece "구현" 0

$d.Bookmarks.Add("__DdeLink__18_641832499", $crouchRng)

# ------------------------------------------------------------------
# 2. Strikethrough "Checking Score"
# ------------------------------------------------------------------
$d.Paragraphs.Item(4).Range.Font.StrikeThrough = 1

# ------------------------------------------------------------------
# 3. Strikethrough "플레이어 움직임 구현" (keeps its red color)
# ------------------------------------------------------------------
$d.Paragraphs.Item(8).Range.Font.StrikeThrough = 1

# ------------------------------------------------------------------
# 4. Strikethrough "히트 체크 구현" (keeps its red color) -- now
#    paragraph 10, since the new paragraph above pushed it down one
#    slot.
# ------------------------------------------------------------------
$d.Paragraphs.Item(10).Range.Font.StrikeThrough = 1

Write-Output "done"
